{"js": "// Replace the three-digit x one-digit multiplication problems/answers\n// in the table with a new generated set of problems, preserving\n// formatting (font, size, alignment) since we replace text in place.\nconst replacements = [\n  { from: \"924\u00d79=8316\", to: \"621\u00d79=5589\" },\n  { from: \"425\u00d75=2125\", to: \"955\u00d77=6685\" },\n  { from: \"177\u00d73=531\", to: \"297\u00d72=594\" },\n  { from: \"549\u00d78=4392\", to: \"692\u00d73=2076\" },\n  { from: \"996\u00d76=5976\", to: \"253\u00d72=506\" },\n  { from: \"647\u00d76=3882\", to: \"263\u00d78=2104\" },\n  { from: \"544\u00d77=3808\", to: \"809\u00d79=7281\" },\n  { from: \"774\u00d79=6966\", to: \"277\u00d76=1662\" },\n  { from: \"189\u00d76=1134\", to: \"291\u00d78=2328\" },\n  { from: \"922\u00d72=1844\", to: \"295\u00d75=1475\" },\n  { from: \"796\u00d74=3184\", to: \"323\u00d76=1938\" },\n  { from: \"873\u00d75=4365\", to: \"784\u00d75=3920\" },\n  { from: \"780\u00d74=3120\", to: \"736\u00d77=5152\" },\n  { from: \"661\u00d73=1983\", to: \"623\u00d79=5607\" },\n  { from: \"803\u00d79=7227\", to: \"687\u00d73=2061\" },\n  { from: \"555\u00d75=2775\", to: \"342\u00d72=684\" },\n  { from: \"722\u00d74=2888\", to: \"781\u00d78=6248\" },\n  { from: \"483\u00d79=4347\", to: \"979\u00d77=6853\" },\n  { from: \"797\u00d73=2391\", to: \"106\u00d75=530\" },\n  { from: \"371\u00d79=3339\", to: \"387\u00d76=2322\" },\n  { from: \"753\u00d77=5271\", to: \"352\u00d75=1760\" },\n  { from: \"778\u00d72=1556\", to: \"526\u00d76=3156\" },\n  { from: \"379\u00d73=1137\", to: \"502\u00d74=2008\" },\n  { from: \"182\u00d72=364\", to: \"622\u00d74=2488\" },\n  { from: \"337\u00d78=2696\", to: \"703\u00d79=6327\" },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems/answers\n# in the table with a new generated set of problems, preserving\n# formatting (font, size, alignment) since we replace text in place.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Find = \"924\u00d79=8316\"; Replace = \"621\u00d79=5589\" },\n  @{ Find = \"425\u00d75=2125\"; Replace = \"955\u00d77=6685\" },\n  @{ Find = \"177\u00d73=531\"; Replace = \"297\u00d72=594\" },\n  @{ Find = \"549\u00d78=4392\"; Replace = \"692\u00d73=2076\" },\n  @{ Find = \"996\u00d76=5976\"; Replace = \"253\u00d72=506\" },\n  @{ Find = \"647\u00d76=3882\"; Replace = \"263\u00d78=2104\" },\n  @{ Find = \"544\u00d77=3808\"; Replace = \"809\u00d79=7281\" },\n  @{ Find = \"774\u00d79=6966\"; Replace = \"277\u00d76=1662\" },\n  @{ Find = \"189\u00d76=1134\"; Replace = \"291\u00d78=2328\" },\n  @{ Find = \"922\u00d72=1844\"; Replace = \"295\u00d75=1475\" },\n  @{ Find = \"796\u00d74=3184\"; Replace = \"323\u00d76=1938\" },\n  @{ Find = \"873\u00d75=4365\"; Replace = \"784\u00d75=3920\" },\n  @{ Find = \"780\u00d74=3120\"; Replace = \"736\u00d77=5152\" },\n  @{ Find = \"661\u00d73=1983\"; Replace = \"623\u00d79=5607\" },\n  @{ Find = \"803\u00d79=7227\"; Replace = \"687\u00d73=2061\" },\n  @{ Find = \"555\u00d75=2775\"; Replace = \"342\u00d72=684\" },\n  @{ Find = \"722\u00d74=2888\"; Replace = \"781\u00d78=6248\" },\n  @{ Find = \"483\u00d79=4347\"; Replace = \"979\u00d77=6853\" },\n  @{ Find = \"797\u00d73=2391\"; Replace = \"106\u00d75=530\" },\n  @{ Find = \"371\u00d79=3339\"; Replace = \"387\u00d76=2322\" },\n  @{ Find = \"753\u00d77=5271\"; Replace = \"352\u00d75=1760\" },\n  @{ Find = \"778\u00d72=1556\"; Replace = \"526\u00d76=3156\" },\n  @{ Find = \"379\u00d73=1137\"; Replace = \"502\u00d74=2008\" },\n  @{ Find = \"182\u00d72=364\"; Replace = \"622\u00d74=2488\" },\n  @{ Find = \"337\u00d78=2696\"; Replace = \"703\u00d79=6327\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $true, $true, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
